$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I ("Quantity") to hold the new
# "Option Type" field.
$ws.Columns.Item(9).Insert()

# Match the width of the neighbouring "Quantity" column that got pushed
# to column J (was column I, 15.68 characters wide).
$ws.Columns.Item(9).ColumnWidth = 14.84

# Header for the new column.
$ws.Cells.Item(1, 9).Value = "Option Type"

# Populate the option type for the two "Options" instrument rows.
$ws.Cells.Item(6, 9).Value = "Regular"
$ws.Cells.Item(7, 9).Value = "Phantom"

$ws.Range("I7").Select()
